# Generate Report for Handoff
# The file "7e7af4e9-8d4e-4713-bf37-157460121337.md" has finished translation
# and is now ready for handoff. Update its status and handoff timestamps
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-14-12 16:14:58"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-12 16:14:55"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-12 16:14:58"
